$d = $word.ActiveDocument

# Merge multi-run paragraphs into single runs by re-finding & replacing each paragraph's full text with itself.
# (Mirrors the original edit session touching each paragraph, which collapses Word's internal run fragmentation
#  and clears stale proofing-error markers.)
$d.Content.Find.Execute(">> Judging from the database shown by the given diagram; it can be seen here that SQL; a special-purpose programming language designed to handle data in a relational database management system has remained consistent as the leading demand of IT skill in throughout these particular time period. Moreover, a text-based programming language such as JavaScript additionally allows both the client-side and sever-side to make web pages interactive. With all of these specific IT skills mentioned; fortunately, all the personally favored jobs above fits in perfectly with these requirements.", $false, $false, $false, $false, $false, $true, 1, $false, ">> Judging from the database shown by the given diagram; it can be seen here that SQL; a special-purpose programming language designed to handle data in a relational database management system has remained consistent as the leading demand of IT skill in throughout these particular time period. Moreover, a text-based programming language such as JavaScript additionally allows both the client-side and sever-side to make web pages interactive. With all of these specific IT skills mentioned; fortunately, all the personally favored jobs above fits in perfectly with these requirements.", 2) | Out-Null
$d.Content.Find.Execute(">> From the information displayed and provided by Burning Glass, without a doubt; Communication skills has won the favoritism in being the number required baseline skills set along with Problem Solving, Organizational Skills, Writing and Teamwork/Collaboration. All of which slots in cleanly with the demanded skills in the listed occupations above, further solidifying the importance of improving said skills. ", $false, $false, $false, $false, $false, $true, 1, $false, ">> From the information displayed and provided by Burning Glass, without a doubt; Communication skills has won the favoritism in being the number required baseline skills set along with Problem Solving, Organizational Skills, Writing and Teamwork/Collaboration. All of which slots in cleanly with the demanded skills in the listed occupations above, further solidifying the importance of improving said skills. ", 2) | Out-Null
$d.Content.Find.Execute(">> Microsoft Windows: broadly known as Windows, is a group of several proprietary graphical operating system families, all of which are developed and marketed by Microsoft. It can be said here that, all the three aforementioned career paths above do not necessarily revolve around this system as they tend to be more obligated to be focused on other aspects of their jobs such as coding, analytical problem-solving, researching and more.", $false, $false, $false, $false, $false, $true, 1, $false, ">> Microsoft Windows: broadly known as Windows, is a group of several proprietary graphical operating system families, all of which are developed and marketed by Microsoft. It can be said here that, all the three aforementioned career paths above do not necessarily revolve around this system as they tend to be more obligated to be focused on other aspects of their jobs such as coding, analytical problem-solving, researching and more.", 2) | Out-Null
$d.Content.Find.Execute("Microsoft C#: is an object-oriented programming language from Microsoft that aims to achieve combination of the computing power of C++ and the programming ease of Visual Basic. Likewise, It is an absence feature that all the occupations written above shared in common and held a fewer or lesser importance. ", $false, $false, $false, $false, $false, $true, 1, $false, "Microsoft C#: is an object-oriented programming language from Microsoft that aims to achieve combination of the computing power of C++ and the programming ease of Visual Basic. Likewise, It is an absence feature that all the occupations written above shared in common and held a fewer or lesser importance. ", 2) | Out-Null
$d.Content.Find.Execute("Oracle: being a multi-model database management system produced and marketed by Oracle, it is commonly used for running online transaction processing (OLTP), data warehousing (DW) and mixed (OLTP & DW) database workloads. Similar to the other two factors; there are other alternatives preferred or used instead of this, for instance Java along with SQL.", $false, $false, $false, $false, $false, $true, 1, $false, "Oracle: being a multi-model database management system produced and marketed by Oracle, it is commonly used for running online transaction processing (OLTP), data warehousing (DW) and mixed (OLTP & DW) database workloads. Similar to the other two factors; there are other alternatives preferred or used instead of this, for instance Java along with SQL.", 2) | Out-Null
$d.Content.Find.Execute(">> Mentoring: In normal circumstances; people who have chosen and walked the career path of IT would not necessarily need further teaching, as most of the important knowledge would have been learnt and taught during their university years. Not to mention, self-study and research can be made easily done individually regarding any possible lingering questions, therefore it is safe to say that this skill isn’t quite a necessity as it unlike for collaboration and communication skills. ", $false, $false, $false, $false, $false, $true, 1, $false, ">> Mentoring: In normal circumstances; people who have chosen and walked the career path of IT would not necessarily need further teaching, as most of the important knowledge would have been learnt and taught during their university years. Not to mention, self-study and research can be made easily done individually regarding any possible lingering questions, therefore it is safe to say that this skill isn’t quite a necessity as it unlike for collaboration and communication skills. ", 2) | Out-Null
$d.Content.Find.Execute("Presentation Skills: Similar with the previously mentioned requirements, this one is not particularly needed either as most projects are either done in teams which only need in-group discussion or workload that requires only one person. ", $false, $false, $false, $false, $false, $true, 1, $false, "Presentation Skills: Similar with the previously mentioned requirements, this one is not particularly needed either as most projects are either done in teams which only need in-group discussion or workload that requires only one person. ", 2) | Out-Null
$d.Content.Find.Execute("Management: Regarding the three specific occupations, management does not seem to be the main focal point of the job, rather they each have their own tasks to do that may overlap with analytical thinking together with working in teams but not much to the point of managing the whole process on their own.", $false, $false, $false, $false, $false, $true, 1, $false, "Management: Regarding the three specific occupations, management does not seem to be the main focal point of the job, rather they each have their own tasks to do that may overlap with analytical thinking together with working in teams but not much to the point of managing the whole process on their own.", 2) | Out-Null

# Append the missing closing paragraph at the end of the document.
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$newParaRange = $d.Paragraphs.Last.Range
$newParaRange.Text = ">>Despite after seeing the charts given by Burning Glass, we have simultaneously came into a collaborative agreement in which none of our favored occupations currently in mind have been changed. In our most humble opinions, staying faithful with the decision you have made along with working hard in regards to improving your skills is better than constantly changing your ideal job in order to adapt with the ever-changing demands. This can be applied beyond jobs selection as it can be proven in normal everyday life, as there are plenty of things you can’t control in life; such as the people around you together with what society requires from you. Yet that doesn’t mean you’ll have to follow to every single one of them, the choice you made will pave a way to a future in which you have hoped and planned for. Personally, it’s preferable to live life the way you want to in contrast of the mindset to continuously follow every new trend, therefore it is why our group members have ultimately stayed with our selected career for now; as possible occurrences in the future including family and personal issues may hinder with the chosen path."
